$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 2.55
$ws.Range("H2").Value = 3.6
$ws.Range("I2").Value = 2.6
$ws.Range("J2").Value = 3.1
$ws.Range("K2").Value = 2.3
$ws.Range("L2").Value = 3.1
$ws.Range("M2").Value = 1.04
$ws.Range("N2").Value = 13
$ws.Range("O2").Value = 1.2
$ws.Range("P2").Value = 4.5
$ws.Range("Q2").Value = 1.67
$ws.Range("R2").Value = 2.2
$ws.Range("S2").Value = 1.33
$ws.Range("T2").Value = 3.25
$ws.Range("U2").Value = 1.57
$ws.Range("V2").Value = 2.25
$ws.Range("W2").Value = 12
$ws.Range("X2").Value = 17
$ws.Range("Z2").Value = 29
$ws.Range("AA2").Value = 21
$ws.Range("AB2").Value = 26
$ws.Range("AC2").Value = 15
$ws.Range("AD2").Value = 7.5
$ws.Range("AG2").Value = 151
$ws.Range("AH2").Value = 12
$ws.Range("AI2").Value = 17
$ws.Range("AK2").Value = 29
$ws.Range("AO2").Value = 13
$ws.Range("AP2").Value = 21
$ws.Range("AS2").Value = 126
$ws.Range("AT2").Value = 3.25
$ws.Range("AW2").Value = 401
$ws.Range("AX2").Value = 4.75
$ws.Range("AY2").Value = 13
$ws.Range("AZ2").Value = 21
$ws.Range("BC2").Value = 126

# Row 3
$ws.Range("G3").Value = 2.1
$ws.Range("I3").Value = 3.1
$ws.Range("J3").Value = 2.6
$ws.Range("L3").Value = 3.4
$ws.Range("Q3").Value = 1.44
$ws.Range("R3").Value = 2.75
$ws.Range("S3").Value = 1.25
$ws.Range("T3").Value = 3.75
$ws.Range("U3").Value = 1.4
$ws.Range("V3").Value = 2.75
$ws.Range("W3").Value = 15
$ws.Range("X3").Value = 15
$ws.Range("Y3").Value = 10
$ws.Range("Z3").Value = 23
$ws.Range("AC3").Value = 21
$ws.Range("AD3").Value = 8.5
$ws.Range("AE3").Value = 12
$ws.Range("AI3").Value = 23
$ws.Range("AJ3").Value = 13
$ws.Range("AK3").Value = 41
$ws.Range("AL3").Value = 23
$ws.Range("AM3").Value = 26
$ws.Range("AN3").Value = 4.75
$ws.Range("AO3").Value = 11
$ws.Range("AT3").Value = 3.75
$ws.Range("AU3").Value = 7
$ws.Range("AV3").Value = 41
$ws.Range("AW3").Value = 251
$ws.Range("AY3").Value = 15
$ws.Range("AZ3").Value = 19
$ws.Range("BB3").Value = 51

# Row 4
$ws.Range("G4").Value = 3.5
$ws.Range("I4").Value = 2.1
$ws.Range("J4").Value = 4
$ws.Range("N4").Value = 9.5
$ws.Range("Q4").Value = 2.03
$ws.Range("R4").Value = 1.78
$ws.Range("U4").Value = 1.8
$ws.Range("V4").Value = 1.91
$ws.Range("X4").Value = 17
$ws.Range("AC4").Value = 9.5
$ws.Range("AH4").Value = 7.5
$ws.Range("AI4").Value = 10
$ws.Range("AK4").Value = 19
$ws.Range("AO4").Value = 19
$ws.Range("AR4").Value = 81
$ws.Range("AY4").Value = 12

# Row 5
$ws.Range("AC5").Value = 9
$ws.Range("AH5").Value = 10
$ws.Range("AP5").Value = 23

# Row 9
$ws.Range("M9").Value = 1.08
$ws.Range("N9").Value = 8

# Row 13
$ws.Range("H13").Value = 4.5
$ws.Range("O13").Value = 1.22
$ws.Range("P13").Value = 4
$ws.Range("Q13").Value = 1.75
$ws.Range("R13").Value = 2.05
$ws.Range("U13").Value = 2
$ws.Range("V13").Value = 1.75
$ws.Range("AC13").Value = 12
$ws.Range("AD13").Value = 9
$ws.Range("AK13").Value = 81
$ws.Range("BB13").Value = 151

# Row 14
$ws.Range("G14").Value = 6.5
$ws.Range("H14").Value = 4.33
$ws.Range("I14").Value = 1.48
$ws.Range("L14").Value = 2
$ws.Range("O14").Value = 1.18
$ws.Range("P14").Value = 4.5
$ws.Range("Q14").Value = 1.62
$ws.Range("R14").Value = 2.25
$ws.Range("S14").Value = 1.3
$ws.Range("T14").Value = 3.4
$ws.Range("Y14").Value = 21
$ws.Range("AT14").Value = 3.4
$ws.Range("AY14").Value = 7
$ws.Range("BC14").Value = 101

# Row 17
$ws.Range("N17").Value = 15
$ws.Range("Y17").Value = 11
$ws.Range("AH17").Value = 29
$ws.Range("AI17").Value = 67
$ws.Range("AL17").Value = 101
$ws.Range("AU17").Value = 12

# Row 21
$ws.Range("H21").Value = 5
$ws.Range("K21").Value = 2.4
$ws.Range("Q21").Value = 1.7
$ws.Range("R21").Value = 2.1
$ws.Range("S21").Value = 1.33
$ws.Range("T21").Value = 3.25
$ws.Range("U21").Value = 2.1
$ws.Range("V21").Value = 1.67
$ws.Range("AC21").Value = 12
$ws.Range("AD21").Value = 9.5
$ws.Range("AH21").Value = 6.5
$ws.Range("AI21").Value = 6
$ws.Range("AJ21").Value = 9
$ws.Range("AL21").Value = 12
$ws.Range("AN21").Value = 9.5
$ws.Range("AS21").Value = 401
$ws.Range("AT21").Value = 3.25
$ws.Range("BC21").Value = 151

Write-Output "Applied 137 cell updates"
